$d = $word.ActiveDocument

# 1. Remove the "_GoBack" bookmark that wraps the last (Heading9) paragraph,
#    leaving the paragraph itself (and its Heading9 style) intact but empty.
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# 2. Append two new body paragraphs (default/Normal style) after that
#    paragraph, before the section break, describing the theme fonts.

# -- First new paragraph: "Calibri (Body)"
$lastPara = $d.Paragraphs.Last
[void]$lastPara.Range.InsertParagraphAfter()
$p1 = $d.Paragraphs.Last
[void]$p1.Range.InsertXML(@"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:t>Calibri</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (Body)</w:t></w:r></w:p>
"@)

# -- Second new paragraph: "Calibri Light (Headings)" in the major theme font
$p1 = $d.Paragraphs.Last
[void]$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Last
[void]$p2.Range.InsertXML(@"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/></w:rPr><w:t>Calibri</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/></w:rPr><w:t>Light</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/></w:rPr><w:t xml:space="preserve"> (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/></w:rPr><w:t>Headings</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/></w:rPr><w:t>)</w:t></w:r></w:p>
"@)

Write-Output "Done"
